# The "kind" column (ID / kind / functional_unit) on the "info" sheet is no
# longer needed now that StreamImpactItem accounts for WasteStream impacts
# directly, so drop it: the "info" table becomes just ID / functional_unit.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("info")

$ws.Columns.Item(2).Delete() | Out-Null

# Leave the selection where the author left it after editing.
$ws.Range("E11").Select() | Out-Null
